$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Python")

# --- Row 48: mark the existing "Sum of Mutated Array Closest to Target" row as solved (Y) ---
$ws.Range("A48").Value = "Y"

# --- Row 49: new entry "Best Sightseeing Pair" ---
$ws.Range("A49").Value = "?"
$ws.Range("B49").Value = "Best Sightseeing Pair"
$ws.Range("C49").Value = 1014
$ws.Range("L49").NumberFormat = "@"
$ws.Range("L49").Value = "05/02/2020"

# --- Row 50: new entry "Best Time to Buy and Sell Stock" ---
$ws.Range("A50").Value = "*"
$ws.Range("B50").Value = "Best Time to Buy and Sell Stock"
$ws.Range("C50").Value = 121
$ws.Range("L50").NumberFormat = "@"
$ws.Range("L50").Value = "05/02/2020"

# Move / restore selection like a user tabbing past the last filled cell
$ws.Range("M50").Select()
